$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.097.34"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "3.625.86"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.74"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.73"
$ws.Range("E6").Value = "  +0.79%  "

$ws.Range("D7").Value = "3.621.09"
$ws.Range("E7").Value = "  -1.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.686"
$ws.Range("E10").Value = "  -2.44%  "

$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.12"
$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000296"
$ws.Range("E13").Value = "  +8.39%  "

$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "4.203.13"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "3.621.53"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.60"
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "67.997.05"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.87"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.45"
$ws.Range("E23").Value = "  +21.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  -3.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.33"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.96"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.62"
$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("E28").Value = "  +5.52%  "

$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  +12.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.27"
$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.62"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "686.57"
$ws.Range("E33").Value = "  +11.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.32"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("E35").Value = "  +1.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.98"
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.99"
$ws.Range("E37").Value = "  -4.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.420"
$ws.Range("E38").Value = "  +5.67%  "

$ws.Range("D39").Value = "0.0₃0800"
$ws.Range("E39").Value = "  +3.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  +16.22%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  +7.09%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.185.40"
$ws.Range("E43").Value = "  +13.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.135"
$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0423"
$ws.Range("E46").Value = "  -1.45%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.133"
$ws.Range("E47").Value = "  -2.54%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.88"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.57"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("E50").Value = "  -3.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.57"
$ws.Range("E51").Value = "  +1.02%  "
